$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last "seeder" row (56) down onto the 7 new
# placeholder rows (57-63) so they pick up the same fonts/fills/number
# formats as the existing M1-W0xx rows.
$ws.Range("A56:AD56").Copy() | Out-Null
$ws.Range("A57:AD63").PasteSpecial(-4122) | Out-Null

# New seeder ids for the newly added rows.
$ids = @("M1-W056", "M1-W057", "M1-W058", "M1-W059", "M1-W060", "M1-W061", "M1-W062")

for ($i = 0; $i -lt 7; $i++) {
    $r = 57 + $i
    $ws.Range("A$r").Value = $ids[$i]

    $dataRow = $ws.Range("B" + $r + ":AD" + $r)
    $dataRow.Value = 0

    # Give the data cells of the new rows their own thin outer border
    # (no top/bottom line - those come from the thick row separators),
    # matching the other seeder rows.
    $dataRow.Borders.Item(8).LineStyle = -4142
    $dataRow.Borders.Item(9).LineStyle = -4142
}

$ws.Range("A1").Select() | Out-Null
